$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 42790
$ws.Range("E2").Value = 838715953858
$ws.Range("F2").Value = 14307408035
$ws.Range("G2").Value = 0.37919

$ws.Range("D3").Value = 2542.26
$ws.Range("E3").Value = 305427899450
$ws.Range("F3").Value = 19893896916
$ws.Range("G3").Value = -0.10441

$ws.Range("E4").Value = 95035237385
$ws.Range("F4").Value = 27751424656
$ws.Range("G4").Value = -0.03898

$ws.Range("D5").Value = 303.96
$ws.Range("E5").Value = 46804906731
$ws.Range("F5").Value = 563647300
$ws.Range("G5").Value = 1.43925

$ws.Range("D6").Value = 96.92
$ws.Range("E6").Value = 41894712720
$ws.Range("F6").Value = 2011983791
$ws.Range("G6").Value = 5.24875

$ws.Range("D7").Value = 0.5768760000000001
$ws.Range("E7").Value = 31268733409
$ws.Range("F7").Value = 507653248
$ws.Range("G7").Value = 0.2911

$ws.Range("D8").Value = 1.001
$ws.Range("E8").Value = 25477897933
$ws.Range("F8").Value = 9250266594
$ws.Range("G8").Value = -0.01051

$ws.Range("D9").Value = 2539.31
$ws.Range("E9").Value = 23612768623
$ws.Range("F9").Value = 18193396
$ws.Range("G9").Value = -0.11701

$ws.Range("D10").Value = 0.5448809999999999
$ws.Range("E10").Value = 19120029255
$ws.Range("F10").Value = 318653867
$ws.Range("G10").Value = -0.53882

$ws.Range("D11").Value = 36.42
$ws.Range("E11").Value = 13366410616
$ws.Range("F11").Value = 509010540
$ws.Range("G11").Value = 1.19509

$ws.Range("D12").Value = 0.082804
$ws.Range("E12").Value = 11827898157
$ws.Range("F12").Value = 445893724
$ws.Range("G12").Value = 3.23823

$ws.Range("B13").Value = 'TRX'
$ws.Range("C13").Value = 'TRON'
$ws.Range("D13").Value = 0.114738
$ws.Range("E13").Value = 10111048246
$ws.Range("F13").Value = 792259148
$ws.Range("G13").Value = 1.50546

$ws.Range("B14").Value = 'DOT'
$ws.Range("C14").Value = 'Polkadot'
$ws.Range("D14").Value = 7.57
$ws.Range("E14").Value = 9989356401
$ws.Range("F14").Value = 169120617
$ws.Range("G14").Value = -0.95464

$ws.Range("B15").Value = 'LINK'
$ws.Range("C15").Value = 'Chainlink'
$ws.Range("D15").Value = 15.04
$ws.Range("E15").Value = 8556595787
$ws.Range("F15").Value = 552633744
$ws.Range("G15").Value = 5.35931

$ws.Range("B16").Value = 'MATIC'
$ws.Range("C16").Value = 'Polygon'
$ws.Range("D16").Value = 0.867233
$ws.Range("E16").Value = 8041587333
$ws.Range("F16").Value = 315590438
$ws.Range("G16").Value = -0.62422

$ws.Range("D17").Value = 2.28
$ws.Range("E17").Value = 7892987706
$ws.Range("F17").Value = 57363484
$ws.Range("G17").Value = 7.93493

$ws.Range("D18").Value = 42781
$ws.Range("E18").Value = 6763821831
$ws.Range("F18").Value = 250284532
$ws.Range("G18").Value = 0.43331

$ws.Range("B19").Value = 'ICP'
$ws.Range("C19").Value = 'Internet Computer'
$ws.Range("D19").Value = 13.3
$ws.Range("E19").Value = 6059479127
$ws.Range("F19").Value = 233563558
$ws.Range("G19").Value = 3.88334

$ws.Range("B20").Value = 'SHIB'
$ws.Range("C20").Value = 'Shiba Inu'
$ws.Range("D20").Value = 0.00000988
$ws.Range("E20").Value = 5819445309
$ws.Range("F20").Value = 117945167
$ws.Range("G20").Value = 0.6561399999999999

$ws.Range("B21").Value = 'LTC'
$ws.Range("C21").Value = 'Litecoin'
$ws.Range("D21").Value = 71.69
$ws.Range("E21").Value = 5312242497
$ws.Range("F21").Value = 514406407
$ws.Range("G21").Value = 0.05971

$ws.Range("B22").Value = 'DAI'
$ws.Range("C22").Value = 'Dai'
$ws.Range("D22").Value = 0.99977
$ws.Range("E22").Value = 5252841442
$ws.Range("F22").Value = 130073276
$ws.Range("G22").Value = -0.02637

$ws.Range("B23").Value = 'BCH'
$ws.Range("C23").Value = 'Bitcoin Cash'
$ws.Range("D23").Value = 255.96
$ws.Range("E23").Value = 5033798569
$ws.Range("F23").Value = 215009984
$ws.Range("G23").Value = -0.13258

$ws.Range("B24").Value = 'UNI'
$ws.Range("C24").Value = 'Uniswap'
$ws.Range("D24").Value = 6.58
$ws.Range("E24").Value = 4952303220
$ws.Range("F24").Value = 159412431
$ws.Range("G24").Value = 0.38575

$ws.Range("B25").Value = 'ETC'
$ws.Range("C25").Value = 'Ethereum Classic'
$ws.Range("D25").Value = 28.05
$ws.Range("E25").Value = 4024338787
$ws.Range("F25").Value = 314351946
$ws.Range("G25").Value = -3.78961

$ws.Range("D26").Value = 10.16
$ws.Range("E26").Value = 3888245163
$ws.Range("F26").Value = 179077629
$ws.Range("G26").Value = 1.09964

$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'LEO Token'
$ws.Range("D27").Value = 4.06
$ws.Range("E27").Value = 3758023659
$ws.Range("F27").Value = 1165327
$ws.Range("G27").Value = 0.10898

$ws.Range("B28").Value = 'OP'
$ws.Range("C28").Value = 'Optimism'
$ws.Range("D28").Value = 3.77
$ws.Range("E28").Value = 3606139832
$ws.Range("F28").Value = 185597732
$ws.Range("G28").Value = -1.15445

$ws.Range("B29").Value = 'NEAR'
$ws.Range("C29").Value = 'NEAR Protocol'
$ws.Range("D29").Value = 3.41
$ws.Range("E29").Value = 3459651464
$ws.Range("F29").Value = 139759161
$ws.Range("G29").Value = 0.66249

$ws.Range("B30").Value = 'XLM'
$ws.Range("C30").Value = 'Stellar'
$ws.Range("D30").Value = 0.11992
$ws.Range("E30").Value = 3394946653
$ws.Range("F30").Value = 56514269
$ws.Range("G30").Value = 0.41485

$ws.Range("B31").Value = 'APT'
$ws.Range("C31").Value = 'Aptos'
$ws.Range("D31").Value = 9.73
$ws.Range("E31").Value = 3244770944
$ws.Range("F31").Value = 317224748
$ws.Range("G31").Value = 0.44356

$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'OKB'
$ws.Range("D32").Value = 53.73
$ws.Range("E32").Value = 3224876054
$ws.Range("F32").Value = 4230407
$ws.Range("G32").Value = 0.2525

$ws.Range("B33").Value = 'INJ'
$ws.Range("C33").Value = 'Injective'
$ws.Range("D33").Value = 38.01
$ws.Range("E33").Value = 3221424410
$ws.Range("F33").Value = 195260724
$ws.Range("G33").Value = 3.00781

$ws.Range("B34").Value = 'TIA'
$ws.Range("C34").Value = 'Celestia'
$ws.Range("D34").Value = 19.72
$ws.Range("E34").Value = 3118861228
$ws.Range("F34").Value = 437259675
$ws.Range("G34").Value = 15.74921

$ws.Range("B35").Value = 'FIL'
$ws.Range("C35").Value = 'Filecoin'
$ws.Range("D35").Value = 6.12
$ws.Range("E35").Value = 3041336526
$ws.Range("F35").Value = 217254639
$ws.Range("G35").Value = 2.90776

$ws.Range("D36").Value = 3.31
$ws.Range("E36").Value = 2938671711
$ws.Range("F36").Value = 109028978
$ws.Range("G36").Value = -2.2742

$ws.Range("B37").Value = 'XMR'
$ws.Range("C37").Value = 'Monero'
$ws.Range("D37").Value = 157.55
$ws.Range("E37").Value = 2857647497
$ws.Range("F37").Value = 84625379
$ws.Range("G37").Value = 3.57223

$ws.Range("D38").Value = 2.07
$ws.Range("E38").Value = 2731666421
$ws.Range("F38").Value = 94796208
$ws.Range("G38").Value = -2.24083

$ws.Range("B39").Value = 'ARB'
$ws.Range("C39").Value = 'Arbitrum'
$ws.Range("D39").Value = 2.14
$ws.Range("E39").Value = 2723260076
$ws.Range("F39").Value = 524765215
$ws.Range("G39").Value = -0.79235

$ws.Range("B40").Value = 'HBAR'
$ws.Range("C40").Value = 'Hedera'
$ws.Range("D40").Value = 0.079858
$ws.Range("E40").Value = 2688359204
$ws.Range("F40").Value = 52655919
$ws.Range("G40").Value = 0.75824

$ws.Range("B41").Value = 'KAS'
$ws.Range("C41").Value = 'Kaspa'
$ws.Range("D41").Value = 0.115743
$ws.Range("E41").Value = 2581956515
$ws.Range("F41").Value = 15714773
$ws.Range("G41").Value = 1.34216

$ws.Range("B42").Value = 'STX'
$ws.Range("C42").Value = 'Stacks'
$ws.Range("D42").Value = 1.67
$ws.Range("E42").Value = 2407120272
$ws.Range("F42").Value = 49566163
$ws.Range("G42").Value = 0.21042

$ws.Range("B43").Value = 'CRO'
$ws.Range("C43").Value = 'Cronos'
$ws.Range("D43").Value = 0.089314
$ws.Range("E43").Value = 2367381478
$ws.Range("F43").Value = 7092169
$ws.Range("G43").Value = 1.39985

$ws.Range("B44").Value = 'MNT'
$ws.Range("C44").Value = 'Mantle'
$ws.Range("D44").Value = 0.715853
$ws.Range("E44").Value = 2242606657
$ws.Range("F44").Value = 94923001
$ws.Range("G44").Value = -5.40362

$ws.Range("D45").Value = 0.03051279
$ws.Range("E45").Value = 2216643529
$ws.Range("F45").Value = 38033027
$ws.Range("G45").Value = -1.68413

$ws.Range("B46").Value = 'TUSD'
$ws.Range("C46").Value = 'TrueUSD'
$ws.Range("D46").Value = 0.994175
$ws.Range("E46").Value = 2037314613
$ws.Range("F46").Value = 95586720
$ws.Range("G46").Value = -0.04741

$ws.Range("D47").Value = 0.997501
$ws.Range("E47").Value = 2015202893
$ws.Range("F47").Value = 2011147005
$ws.Range("G47").Value = -0.2988

$ws.Range("B48").Value = 'MKR'
$ws.Range("C48").Value = 'Maker'
$ws.Range("D48").Value = 2088.48
$ws.Range("E48").Value = 1925059006
$ws.Range("F48").Value = 82932685
$ws.Range("G48").Value = 0.23261

$ws.Range("B49").Value = 'QNT'
$ws.Range("C49").Value = 'Quant'
$ws.Range("D49").Value = 118.84
$ws.Range("E49").Value = 1724872662
$ws.Range("F49").Value = 21658082
$ws.Range("G49").Value = -0.4053

$ws.Range("B50").Value = 'BSV'
$ws.Range("C50").Value = 'Bitcoin SV'
$ws.Range("D50").Value = 87.15000000000001
$ws.Range("E50").Value = 1719689492
$ws.Range("F50").Value = 89331112
$ws.Range("G50").Value = 3.15771

$ws.Range("D51").Value = 0.72019
$ws.Range("E51").Value = 1655177338
$ws.Range("F51").Value = 226226717
$ws.Range("G51").Value = 1.88748
